$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Quantity (C) and Total Cost (E) for rows 3-6
$ws.Range("C3").Value = "1"
$ws.Range("E3").Value = "9.25"

$ws.Range("C4").Value = "1"
$ws.Range("E4").Value = "23.75"

$ws.Range("C5").Value = "1"
$ws.Range("E5").Value = "13.00"

$ws.Range("C6").Value = "1"
$ws.Range("E6").Value = "13.00"
